$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" sheet
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column)
$ws.Columns.Item(14).Insert()

# Set the new column's width to match column M (width 11), no bestFit
$ws.Columns.Item(14).ColumnWidth = 11

# Update selection on this sheet
$ws.Range("K16").Select()

# Make "Repayment schedule" the active/selected sheet
$ws.Activate()

# Remove tabSelected from "NewLoanInput" sheet by activating this one instead
$newLoanInput = $wb.Worksheets.Item("NewLoanInput")
$newLoanInput.Select()
$ws.Select()
